$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they stay literal strings
$numericTextCells = @("D5", "D6", "D11", "D13", "D18", "D19", "D22", "D25", "D27", "D28", "D30", "D34", "D40", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($ref in $numericTextCells) { $ws.Range($ref).NumberFormat = "@" }

# Apply new cell values
$ws.Range("D2").Value = "70.197.46"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "3.939.51"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "610.91"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "171.54"
$ws.Range("E6").Value = "  +6.27%  "
$ws.Range("D7").Value = "3.939.17"
$ws.Range("E7").Value = "  +2.65%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "6.43"
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("D13").Value = "38.70"
$ws.Range("E13").Value = "  +5.30%  "
$ws.Range("E14").Value = "  +6.05%  "
$ws.Range("D15").Value = "4.602.01"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "3.959.64"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "70.215.83"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "7.67"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "18.67"
$ws.Range("E19").Value = "  +9.02%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").Value = "496.10"
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("E23").Value = "  +4.30%  "
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").Value = "86.17"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").Value = "12.38"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "3.01"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").Value = "4.091.40"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "32.39"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").Value = "3.901.68"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").Value = "3.32"
$ws.Range("E40").Value = "  +12.18%  "
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +7.47%  "
$ws.Range("D44").Value = "439.59"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "48.40"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "8.70"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "0.000277"
$ws.Range("E48").Value = "  +22.98%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0368"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "40.78"
$ws.Range("E50").Value = "  +6.01%  "
$ws.Range("D51").Value = "143.47"
$ws.Range("E51").Value = "  +0.12%  "

# Restore default (General) formatting so styles match the original workbook
foreach ($ref in $numericTextCells) { $ws.Range($ref).ClearFormats() }
